$d = $word.ActiveDocument

# Replace "<id>p003r_a3</id>" with "<id>p003r_3</id>" - collapses the
# three runs (<id>, p003r_a3, </id>) into a single run carrying the
# formatting of the first run (Courier New / 7f6000 / sz18).
$d.Content.Find.Execute("<id>p003r_a3</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p003r_3</id>", 2)

# Replace "<id>p003v_a1</id>" with "<id>p003v_1</id>" - same collapse.
$d.Content.Find.Execute("<id>p003v_a1</id>", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "<id>p003v_1</id>", 2)
